# #9 treat x columns as notes
# Add note cells in columns H/I on the "30.12.20" sheet and update the
# active selection to reflect the new focus on the notes column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("30.12.20")

$ws.Range("H2").Value = "blabla"
$ws.Range("I2").Value = "another note"
$ws.Range("I3").Value = "bar"
$ws.Range("H14").Value = "test"

$ws.Activate()
$ws.Range("H14").Select()
